$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (2514) down
# across the full range of new rows (2515:2611) so the new rows pick up
# the same date/currency formatting, fills and borders as prior rows.
$ws.Range("A2514:B2514").Copy()
$ws.Range("A2515:B2611").PasteSpecial(-4122)

# New SEAGRI/CEPEA price rows: date serial (col A) + price (col B)
$newData = @(@(45338,102.13),@(45341,102.33),@(45342,102.27),@(45343,103.76),@(45344,102.25),@(45345,100.33),@(45348,103),@(45349,108.38),@(45350,101.63),@(45351,100.5),@(45352,99.75),@(45355,100.38),@(45356,101.38),@(45357,101.63),@(45358,103),@(45359,102.75),@(45362,104.63),@(45363,105),@(45364,107),@(45365,107.8),@(45366,107.8),@(45369,107),@(45370,106.5),@(45371,106.5),@(45372,106.5),@(45373,107.75),@(45376,107.88),@(45377,108.25),@(45378,107.88),@(45379,106.63),@(45383,106.63),@(45384,107.88),@(45385,108.75),@(45386,108.63),@(45387,108),@(45390,108.75),@(45391,110.38),@(45392,108.67),@(45393,109.88),@(45394,109.75),@(45397,110.5),@(45398,110.5),@(45399,112.25),@(45400,112.75),@(45401,112),@(45404,111.25),@(45405,111.25),@(45406,112),@(45407,112.25),@(45408,112.13),@(45411,111.63),@(45412,109.75),@(45414,110.75),@(45415,111.5),@(45418,111.38),@(45419,111.38),@(45420,115.5),@(45421,115.38),@(45422,114.88),@(45425,114.75),@(45426,115.63),@(45427,115.63),@(45428,115.63),@(45429,115.63),@(45432,115.63),@(45433,115.63),@(45434,117.25),@(45435,117.5),@(45436,121),@(45439,119.88),@(45440,118.47),@(45441,120.38),@(45446,120.38),@(45447,120.38),@(45448,120.38),@(45449,120.38),@(45450,120.38),@(45451,120.38),@(45453,119.53),@(45454,119.53),@(45455,119.53),@(45456,119.75),@(45457,120.25),@(45460,119.63),@(45460,119.63),@(45461,119.63),@(45462,120.88),@(45463,121.88),@(45468,121),@(45469,120),@(45470,121.38),@(45471,120.75),@(45474,121.5),@(45476,121.38),@(45477,122.5),@(45478,122.5),@(45481,122.5))

$startRow = 2515
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
}

# Match the row heights of the new rows to the existing data rows
$ws.Range("A2515:A2611").RowHeight = 15

# Restore the active selection / scroll position saved with the workbook
$ws.Application.Goto($ws.Range("B2610"))
$ws.Range("B2610").Select()

Write-Output "Added $($newData.Count) rows (2515-2611)"
